$d = $word.ActiveDocument

# The paragraph's only run currently holds an inline picture (InlineShape).
# Replace it with the literal text "/getid."
$shape = $d.InlineShapes.Item(1)
$range = $shape.Range
$range.Text = "/getid."
